# Update countries & provincias Spain
# - Update the "last updated" timestamp (22:50 -> 23:20)
# - Refresh a number of per-country stat rows with new counts
# - Re-rank three countries (Rumania, Ucrania, Nigeria) whose totals grew
#   past their neighbours, so they move up in the (descending, sorted by
#   "Casos totales") table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Timestamp banner (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 23:20"

# ---------------------------------------------------------------------
# 2) Plain stat refreshes (country keeps its rank / row position)
# ---------------------------------------------------------------------
# Canada (row 18)
$ws.Cells.Item(18,2).Value = 7427
$ws.Cells.Item(18,3).Value = 1107
$ws.Cells.Item(18,4).Value = 1093
$ws.Cells.Item(18,5).Value = 6250
$ws.Cells.Item(18,7).Value = 19
$ws.Cells.Item(18,8).Value = 84

# India (row 43)
$ws.Cells.Item(43,4).Value = 102
$ws.Cells.Item(43,5).Value = 1117

# Costa Rica (row 77)
$ws.Cells.Item(77,2).Value = 330
$ws.Cells.Item(77,3).Value = 16
$ws.Cells.Item(77,4).Value = 4
$ws.Cells.Item(77,5).Value = 324
$ws.Cells.Item(77,6).Value = 7

# Venezuela (row 105)
$ws.Cells.Item(105,2).Value = 135
$ws.Cells.Item(105,3).Value = 16
$ws.Cells.Item(105,5).Value = 93

# ---------------------------------------------------------------------
# 3) Re-ranked countries: insert a fresh row higher up with the new
#    totals, then delete the now-duplicate row that was pushed down.
# ---------------------------------------------------------------------

# --- Rumania: was row 33 (between Ecuador and Japon), moves to row 30
#     (right after Chile, before Polonia) ---
$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30,1).Value = "Rumania"
$ws.Cells.Item(30,2).Value = 2109
$ws.Cells.Item(30,3).Value = 294
$ws.Cells.Item(30,4).Value = 209
$ws.Cells.Item(30,5).Value = 1835
$ws.Cells.Item(30,6).Value = 33
$ws.Cells.Item(30,7).Value = 22
$ws.Cells.Item(30,8).Value = 65
$ws.Rows.Item(34).Delete()

# --- Ucrania: was row 69 (between Armenia and Hungria), moves to row 65
#     (right after Argelia, before Marruecos) ---
$ws.Rows.Item(65).Insert()
$ws.Cells.Item(65,1).Value = "Ucrania"
$ws.Cells.Item(65,2).Value = 548
$ws.Cells.Item(65,3).Value = 73
$ws.Cells.Item(65,4).Value = 8
$ws.Cells.Item(65,5).Value = 527
$ws.Cells.Item(65,6).Value = 0
$ws.Cells.Item(65,7).Value = 3
$ws.Cells.Item(65,8).Value = 13
$ws.Rows.Item(70).Delete()

# --- Nigeria: was row 110 (between Estado de Palestina and Camboya),
#     moves to row 106 (right after Venezuela, before Mauricio) ---
$ws.Rows.Item(106).Insert()
$ws.Cells.Item(106,1).Value = "Nigeria"
$ws.Cells.Item(106,2).Value = 131
$ws.Cells.Item(106,3).Value = 20
$ws.Cells.Item(106,4).Value = 8
$ws.Cells.Item(106,5).Value = 121
$ws.Cells.Item(106,6).Value = 0
$ws.Cells.Item(106,7).Value = 1
$ws.Cells.Item(106,8).Value = 2
$ws.Rows.Item(111).Delete()
